# Generate Report for Handback
# Updates the "Latest Handback DateTime" (column K) for the
# b64778cb-8792-415b-80c6-c326caff7005.md row on both the zh-cn and
# de-de localization-status sheets, reflecting a newly generated
# handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-09-07 03:28:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-07 03:28:23"
